# "tabelas atualizadas com novos dados" -- update the results table with the
# refreshed counts / percentages / p-values.
#
# The first (and only) table in the document holds the data; each changed
# value is addressed by its (row, column) position so the edit is
# unambiguous regardless of any other formatting in the doc.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-Cell($row, $col, $value) {
    $t.Cell($row, $col).Range.Text = $value
}

# Row 2 ("n"): F column
Set-Cell 2 2 "263"

# Row 3 ("Profilaxia (%)"): F column, p-value
Set-Cell 3 2 "238 (94.8)"
Set-Cell 3 4 "0.818"

# Row 4 ("Dabigatrana (%)"): p-value
Set-Cell 4 4 "0.376"

# Row 5 ("Enoxaparina (%)"): F column, p-value
Set-Cell 5 2 "236 (94.0)"
Set-Cell 5 4 "0.422"

# Row 6 ("Rivoraxabana (%)"): F column, p-value
Set-Cell 6 2 "57 (22.5)"
Set-Cell 6 4 "0.322"

# Row 7 ("Warfarina (%)"): F column, p-value
Set-Cell 7 2 "115 (50.9)"
Set-Cell 7 4 "0.751"
